# Changes from sprint 1 review
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "GUI interface" backlog item (row 7, Item # 4) no longer has a
# description, estimate, or "by" assignment - clear those cells.
$ws.Range("C7").Value = $null
$ws.Range("D7").Value = $null
$ws.Range("E7").Value = $null

# Leave a couple of sprint-review notes off to the side, in column G.
$ws.Range("G7").Value = "User stories only"
$ws.Range("G8").Value = "functional"

# Leave the selection on the last cell that was touched.
$ws.Range("G9").Select()
